$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -3.5056
$ws.Range("C3").Value = 34.4601
$ws.Range("C4").Value = -113.9284
$ws.Range("C5").Value = -81.8764
$ws.Range("C6").Value = -68.895
$ws.Range("C7").Value = -90.8613
$ws.Range("C8").Value = -263.1822
$ws.Range("C9").Value = -391.5661
$ws.Range("C10").Value = -60.9285
$ws.Range("C11").Value = 247.0124
$ws.Range("C12").Value = 27.7723
$ws.Range("C13").Value = -80.4056
$ws.Range("C14").Value = 266.2759
$ws.Range("C15").Value = 174.5275
$ws.Range("C16").Value = 132.9962
$ws.Range("C17").Value = -38.522
$ws.Range("C18").Value = 56.8181
